# "Modelo de negocios hidrologicas" - add a new bibliography entry to row 112
# and extend the log with 25 more blank numbered rows (115-139), matching the
# formatting already used by the existing blank rows (112-114).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- 1. Fill in the new reference in row 112 (Donde / Tipo / Tema) ---
$ws.Range("C112").Value = "https://www.murfreesborotn.gov/1583/Journey-to-The-Tap-How-Water-Gets-to-You"
$ws.Range("E112").Value = "jpg"
$ws.Range("F112").Value = "proceso agua"
$ws.Rows.Item(112).RowHeight = 45.75

# --- 2. Append 25 new blank rows (115-139), cloning the look of row 114 ---
$ws.Range("A114:H114").Copy()
$ws.Range("A115:H139").PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($r = 115; $r -le 139; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
    $ws.Rows.Item($r).RowHeight = 15.75
}

# --- 3. Restore the on-screen view so the new rows are visible/selected ---
$ws.Activate()
$ws.Range("A112:A139").Select()
$excel.ActiveWindow.ScrollRow = 112
